$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 3729
$ws.Cells.Item(5, 6).Value = 3729
$ws.Cells.Item(7, 6).Value = 5265
$ws.Cells.Item(8, 6).Value = 578
$ws.Cells.Item(9, 6).Value = 410
$ws.Cells.Item(10, 6).Value = 219
$ws.Cells.Item(11, 6).Value = 1042
$ws.Cells.Item(13, 6).Value = 128
$ws.Cells.Item(14, 6).Value = 46
$ws.Cells.Item(15, 6).Value = 724
$ws.Cells.Item(16, 6).Value = 353
$ws.Cells.Item(17, 6).Value = 44
$ws.Cells.Item(21, 6).Value = 370
$ws.Cells.Item(22, 6).Value = 6025
$ws.Cells.Item(26, 6).Value = 6599
$ws.Cells.Item(28, 6).Value = 22
$ws.Cells.Item(29, 6).Value = 3252
$ws.Cells.Item(31, 6).Value = 744
$ws.Cells.Item(32, 6).Value = 4455
$ws.Cells.Item(35, 6).Value = 149
$ws.Cells.Item(36, 6).Value = 1120
$ws.Cells.Item(37, 6).Value = 99
$ws.Cells.Item(39, 6).Value = 6
$ws.Cells.Item(40, 6).Value = 913
$ws.Cells.Item(41, 6).Value = 1108
$ws.Cells.Item(42, 6).Value = 2054
$ws.Cells.Item(43, 6).Value = 4

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 64

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 1147
$ws.Cells.Item(4, 6).Value = 53

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 1147
$ws.Cells.Item(5, 6).Value = 53
$ws.Cells.Item(7, 6).Value = 3729
$ws.Cells.Item(8, 6).Value = 3729
$ws.Cells.Item(10, 6).Value = 5265
$ws.Cells.Item(11, 6).Value = 578
$ws.Cells.Item(12, 6).Value = 410
$ws.Cells.Item(13, 6).Value = 219
$ws.Cells.Item(14, 6).Value = 1042
$ws.Cells.Item(16, 6).Value = 128
$ws.Cells.Item(17, 6).Value = 46
$ws.Cells.Item(18, 6).Value = 724
$ws.Cells.Item(19, 6).Value = 353
$ws.Cells.Item(20, 6).Value = 44
$ws.Cells.Item(25, 6).Value = 370
$ws.Cells.Item(26, 6).Value = 6025
$ws.Cells.Item(30, 6).Value = 6608
$ws.Cells.Item(32, 6).Value = 22
$ws.Cells.Item(33, 6).Value = 3252
$ws.Cells.Item(35, 6).Value = 744
$ws.Cells.Item(36, 6).Value = 4455
$ws.Cells.Item(40, 6).Value = 149
$ws.Cells.Item(41, 6).Value = 1120
$ws.Cells.Item(42, 6).Value = 99
$ws.Cells.Item(44, 6).Value = 6
$ws.Cells.Item(45, 6).Value = 913
$ws.Cells.Item(46, 6).Value = 1108
$ws.Cells.Item(48, 6).Value = 2054
$ws.Cells.Item(49, 6).Value = 4
$ws.Cells.Item(50, 6).Value = 64
